# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Update the FullTruckUsage report:
#  - A6 date changes from 45295 to 45294 (2024-01-04 -> 2024-01-03)
#  - E8, E10, E12 truck_used counts change from 1 to 2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FullTruckUsage")

# A6: date serial 45295 -> 45294
$ws.Range("A6").Value = 45294

# E8, E10, E12: truck_used 1 -> 2
$ws.Range("E8").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("E12").Value = 2
